# redesign sso support AD
# Rework the ldap_config staging sheet: replace the numeric "licensee_id"
# column with a textual "domain" column, swap the search_scope /
# account_password header order, restore admin_account / admin_password
# to their intended columns, and populate admin credentials for the
# test_ldap (row 3) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("F1").Value = "account_password"
$ws.Range("J1").Value = "search_scope"
$ws.Range("K1").Value = "admin_account"
$ws.Range("L1").Value = "admin_password"

# --- Row 3 (test_ldap): add the new admin credentials ---
$ws.Range("K3").Value = "ldap.admin@dctest.local"
$ws.Range("K3").Font.Size = 11
$ws.Range("L3").Value = "Cc123456"

# --- Column A: licensee_id -> domain (header + both data rows) ---
$ws.Range("A1").Value = "domain"
$ws.Range("A2").Value = "mo.laxino.com"
$ws.Range("A3").Value = "dctest.local"

# --- Column K width: widen and drop the bestFit auto-size flag ---
$ws.Columns.Item(11).ColumnWidth = 22.1

# --- Final selection left on A3 ---
$ws.Range("A3").Select()
